$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = 4
$ws.Range("B7").Value = "Mahesh"

$ws.Range("B7").Select()
